$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRF = 9.815818181818182

for ($r = 29; $r -le 40; $r++) {
    $ws.Cells.Item($r, 9).Value = $newRF
}
